$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Details")

# Update the candidate's last name (D5) and first name (D6). These feed
# the "Cahier des charges" and "Prop d'engagement" sheets via formulas
# (Details!D5 / Details!D6), which will recalc automatically.
$ws.Range("D5").Value = "MIGUEL"
$ws.Range("D6").Value = "Luis"

# Row 24 (Type of Permit) no longer needs to be as tall.
$ws.Rows.Item(24).RowHeight = 34

# Move the active selection/view down to D26.
$ws.Activate() | Out-Null
$ws.Range("D26").Select() | Out-Null
